# railway_links.xlsx update
#
# Adds newly surveyed railway links to the "Hoja1" table (A: id_link,
# B: distance, C: gauge):
#   - three new "ancha"  links: 27-51, 35-44, 49-1023
#   - two   new "media"  links: 68-69, 69-70
# and relabels two existing "media" link ids whose numbering scheme changed:
#   - 3-12     -> 3-13
#   - 12-1012  -> 13-1012
#
# The three new "ancha" rows are inserted right after the existing "ancha"
# block (before row 97), so every row from the old 97 down to the old 143
# moves down by 3 (new rows 100-146). That block is shifted with plain
# value writes (bottom row first, so each source row is captured before it
# is overwritten) rather than a structural row insert, which keeps the
# sheet's trailing marker row and overall dimension untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: shift rows 97-143 down to 100-146 (bottom-up) ---------------
$shiftedRows = @(
    @("71-73",   353, "media"),
    @("70-71",    80, "media"),
    @("68-70",   340, "media"),
    @("17-1011", 199, "media"),
    @("17-18",   114, "media"),
    @("12-1012",  60, "media"),
    @("12-18",    66, "media"),
    @("12-13",    60, "media"),
    @("11-1012", 120, "media"),
    @("11-1011", 108, "media"),
    @("11-70",   100, "media"),
    @("3-1003",   45, "media"),
    @("3-12",    250, "media"),
    @("1-1003",   58, "media"),
    @("1-3",    94.8, "media"),
    @("91-92",   169, "angosta"),
    @("89-91",   216, "angosta"),
    @("84-1002",  70, "angosta"),
    @("83-1002",  98, "angosta"),
    @("83-84",    50, "angosta"),
    @("81-1001",  36, "angosta"),
    @("80-1001",  80, "angosta"),
    @("78-1001", 142, "angosta"),
    @("78-83",   166, "angosta"),
    @("78-79",   493, "angosta"),
    @("77-80",   196, "angosta"),
    @("76-77",   453, "angosta"),
    @("75-76",   298, "angosta"),
    @("67-68",   182, "angosta"),
    @("66-78",   262, "angosta"),
    @("66-77",   493, "angosta"),
    @("66-67",   260, "angosta"),
    @("65-1002", 116, "angosta"),
    @("64-65",   148, "angosta"),
    @("63-67",   360, "angosta"),
    @("56-1004", 208, "angosta"),
    @("55-89",   250, "angosta"),
    @("55-65",   221, "angosta"),
    @("55-56",   123, "angosta"),
    @("21-1021", 120, "angosta"),
    @("21-1004", 240, "angosta"),
    @("17-21",   172, "angosta"),
    @("15-17",    95, "angosta"),
    @("14-67",   340, "angosta"),
    @("14-63",   158, "angosta"),
    @("14-15",   225, "angosta"),
    @("1-1021",  185, "angosta")
)

$destRow = 146
foreach ($link in $shiftedRows) {
    $ws.Range("A$destRow").Value = $link[0]
    $ws.Range("B$destRow").Value = $link[1]
    $ws.Range("C$destRow").Value = $link[2]
    $destRow = $destRow - 1
}

# --- Step 2: relabel the two renamed "media" ids at their new rows -------
$ws.Range("A134").Value = "3-13"
$ws.Range("A141").Value = "13-1012"

# --- Step 3: new "ancha" rows 97-98 ---------------------------------------
$ws.Range("A97").Value = "27-51"
$ws.Range("B97").Value = 308
$ws.Range("C97").Value = "ancha"
$ws.Range("A98").Value = "35-44"
$ws.Range("B98").Value = 214
$ws.Range("C98").Value = "ancha"

# --- Step 4: new "media" rows 147-148, appended after the existing data --
$ws.Range("A147").Value = "68-69"
$ws.Range("B147").Value = 245
$ws.Range("C147").Value = "media"
$ws.Range("A148").Value = "69-70"
$ws.Range("B148").Value = 210
$ws.Range("C148").Value = "media"

# --- Step 5: new "ancha" row 99 -------------------------------------------
$ws.Range("A99").Value = "49-1023"
$ws.Range("B99").Value = 170
$ws.Range("C99").Value = "ancha"

# --- Step 6: update the sheet selection to match the authored state ------
[void]$ws.Range("C2:C148").Select()
